$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.552.07'
$ws.Range('E2').Value = '  +0.68%  '

$ws.Range('D3').Value = '2.974.63'
$ws.Range('E3').Value = '  +1.95%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '380.20'
$ws.Range('E5').Value = '  +4.13%  '

$ws.Range('D6').Value = '106.10'
$ws.Range('E6').Value = '  +2.42%  '

$ws.Range('D7').Value = '0.545'
$ws.Range('E7').Value = '  +0.84%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('D9').Value = '0.598'
$ws.Range('E9').Value = '  +1.17%  '

$ws.Range('D10').Value = '37.43'
$ws.Range('E10').Value = '  +1.09%  '

$ws.Range('E11').Value = '  +0.45%  '

$ws.Range('D12').Value = '0.0843'
$ws.Range('E12').Value = '  +1.12%  '

$ws.Range('D13').Value = '18.67'
$ws.Range('E13').Value = '  +0.90%  '

$ws.Range('D14').Value = '3.448.21'
$ws.Range('E14').Value = '  +2.19%  '

$ws.Range('D15').Value = '7.51'
$ws.Range('E15').Value = '  +1.87%  '

$ws.Range('D16').Value = '2.989.88'
$ws.Range('E16').Value = '  +2.63%  '

$ws.Range('D17').Value = '0.968'
$ws.Range('E17').Value = '  +1.54%  '

$ws.Range('D18').Value = '51.569.80'
$ws.Range('E18').Value = '  +0.75%  '

$ws.Range('E19').Value = '  +2.14%  '

$ws.Range('D20').Value = '7.40'
$ws.Range('E20').Value = '  +1.96%  '

$ws.Range('D21').Value = '13.01'
$ws.Range('E21').Value = '  -0.31%  '

$ws.Range('D22').Value = '0.0₃0958'
$ws.Range('E22').Value = '  +1.11%  '

$ws.Range('D23').Value = '69.26'
$ws.Range('E23').Value = '  +1.62%  '

$ws.Range('D24').Value = '263.67'
$ws.Range('E24').Value = '  +1.38%  '

$ws.Range('D25').Value = '2.79'
$ws.Range('E25').Value = '  +3.49%  '

$ws.Range('E26').Value = '  -2.39%  '

$ws.Range('D27').Value = '7.25'
$ws.Range('E27').Value = '  +17.82%  '

$ws.Range('D28').Value = '7.43'
$ws.Range('E28').Value = '  +1.58%  '

$ws.Range('E29').Value = '  -0.02%  '

$ws.Range('D30').Value = '26.03'
$ws.Range('E30').Value = '  +0.21%  '

$ws.Range('E31').Value = '  +2.85%  '

$ws.Range('D32').Value = '9.87'
$ws.Range('E32').Value = '  -0.87%  '

$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = '34.72'
$ws.Range('E33').Value = '  -1.56%  '

$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').Value = '2.10'
$ws.Range('E34').Value = '  -2.13%  '

$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.0463'
$ws.Range('E35').Value = '  +9.33%  '

$ws.Range('D36').Value = '51.32'
$ws.Range('E36').Value = '  +1.27%  '

$ws.Range('E37').Value = '  +0.17%  '

$ws.Range('D38').Value = '3.09'
$ws.Range('E38').Value = '  -1.69%  '

$ws.Range('D39').Value = '17.38'
$ws.Range('E39').Value = '  +2.31%  '

$ws.Range('E40').Value = '  -6.76%  '

$ws.Range('E41').Value = '  -1.05%  '

$ws.Range('D42').Value = '0.116'
$ws.Range('E42').Value = '  +2.37%  '

$ws.Range('D43').Value = '123.44'
$ws.Range('E43').Value = '  +4.70%  '

$ws.Range('E44').Value = '  -2.00%  '

$ws.Range('D45').Value = '2.09'
$ws.Range('E45').Value = '  -1.09%  '

$ws.Range('D46').Value = '0.278'
$ws.Range('E46').Value = '  +17.88%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '2.37'
$ws.Range('E47').Value = '  +4.60%  '

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.048.84'
$ws.Range('E48').Value = '  -0.99%  '

$ws.Range('D49').Value = '3.25'
$ws.Range('E49').Value = '  +1.18%  '

$ws.Range('D50').Value = '0.0351'
$ws.Range('E50').Value = '  +9.80%  '

$ws.Range('D51').Value = '5.18'
$ws.Range('E51').Value = '  +3.04%  '
